$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix schema definition cells ---
# Row 2 ("id"): drop the stray duplicate Japanese description in ColumnName
$ws.Range("B2").ClearContents()

# Row 3 ("ringiNo" -> "ringino"): lower-case column name + correct Type/Null
$ws.Range("A3").Value = "ringino"
$ws.Range("C3").Value = "int(10)"
$ws.Range("E3").Value = "NO"

# Row 4 ("approverLayer" -> "approverlayer") + correct Type
$ws.Range("A4").Value = "approverlayer"
$ws.Range("C4").Value = "int(3)"

# Row 5 ("approverID" -> "approverid")
$ws.Range("A5").Value = "approverid"

# Row 6 ("approveDate" -> "approvedate")
$ws.Range("A6").Value = "approvedate"

# Row 7 ("ringStatus" -> "ringistatus") + correct Type
$ws.Range("A7").Value = "ringistatus"
$ws.Range("C7").Value = "varchar(255)"

# Row 8 ("lastLayerFlg" -> "lastlayerflg")
$ws.Range("A8").Value = "lastlayerflg"

# --- Column widths (best-fit sizing for the renamed/shorter columns) ---
$ws.Columns.Item(1).ColumnWidth = 12.92
$ws.Columns.Item(2).ColumnWidth = 21.17
$ws.Columns.Item(3).ColumnWidth = 10.92

# --- Selection moved from D3 to E3 ---
$ws.Range("E3").Select()
